$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LCOE")

# New row 12: only_PV_supply (rerun scenario, copies formatting from row 11)
$ws.Range("A11:H11").Copy($ws.Range("A12:H12"))
$ws.Range("A12").Value = "only_PV_supply"
$ws.Range("B12").Value = 278.40065220000002
$ws.Range("C12").Value = 1538.9369380000001
$ws.Range("D12").Value = 389849663.30000001
$ws.Range("E12").Value = -100014.07429999999
$ws.Range("F12").Value = 98996.813769999993
$ws.Range("G12").Value = 17908.97133
$ws.Range("H12").Value = 14.09394457

# New row 13: only_PV_supply_PV (rerun scenario, copies formatting from row 12)
$ws.Range("A12:H12").Copy($ws.Range("A13:H13"))
$ws.Range("A13").Value = "only_PV_supply_PV"
$ws.Range("B13").Value = 278.40065220000002
$ws.Range("C13").Value = 1538.9369380000001
$ws.Range("D13").Value = 389849663.30000001
$ws.Range("E13").Value = -100014.07429999999
$ws.Range("F13").Value = 98996.813769999993
$ws.Range("G13").Value = 17908.97133
$ws.Range("H13").Value = 14.09394457

# Restore the previously-selected cell on the sheet
$ws.Range("C16").Select()

# Match the saved window geometry from the author's session
$win = $wb.Windows.Item(1)
$win.Width = 19420
$win.Height = 10300
$win.Left = -110
$win.Top = -110
